$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.985.89"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "2.222.09"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0780"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "2.564.52"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.37%  "
$ws.Range("D17").Value = "2.232.43"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.733"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "39.900.75"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  +5.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0716"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0996"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.07%  "
$ws.Range("D42").Value = "2.085.66"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0272"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("E47").Value = "  -8.76%  "
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "2.435.86"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  +2.20%  "
